# Feature : trend finder
# Inserts a new archetype row (City of Gold Haven / Sanctuary Haven) in the
# middle of the table and appends two more new archetype rows at the end
# (City of Gold Haven base entry + Dreadlord Shadow), matching the upstream
# commit that added new meta entries to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at 35 ("Sanctuary Haven") - pushes the old rows 35-47
#    down to 36-48.
# ---------------------------------------------------------------------------
$ws.Rows("35:35").Insert()

# Reuse the existing border/fill/wrap style from row 34 (plain bordered box,
# same as all the other simple rows) so we don't fork a brand-new style.
$ws.Range("A34:H34").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)
$ws.Rows("35:35").RowHeight = 15.75

$ws.Range("A35").Value = "Sanctuary Haven"
$ws.Range("B35").Value = "Haven"
$ws.Range("C35").Value = "Holy Sanctuary"
$ws.Range("D35").Value = "74xJg"
$ws.Range("E35").Value = "Holy Sanctuary"
$ws.Range("F35").Value = "74xJg"
$ws.Range("G35").Value = "None"
$ws.Range("H35").Value = "None"

# ---------------------------------------------------------------------------
# 2) Append two brand-new rows at the bottom of the table: 49 and 50.
# ---------------------------------------------------------------------------

# Row 49 ("City of Gold Haven") matches the style pattern used on row 46
# (name/class boxed + shaded, identifiers plain, anti-identifier boxed).
$ws.Range("A46:H46").Copy()
$ws.Range("A49:H49").PasteSpecial(-4122)
$ws.Rows("49:49").RowHeight = 15.75
$ws.Range("E49").ClearFormats()

$ws.Range("A49").Value = "City of Gold Haven"
$ws.Range("B49").Value = "Haven"
$ws.Range("C49").Value = "City of Gold"
$ws.Range("D49").Value = "6cQBg"
$ws.Range("E49").Value = "City of Gold"
$ws.Range("F49").Value = "6cQBg"
$ws.Range("G49").Value = "None"
$ws.Range("H49").Value = "None"

# Row 50 ("Dreadlord Shadow") matches the style pattern used on row 48
# (name/class boxed + shaded, identifier + anti-identifier boxed + shaded).
$ws.Range("A48:H48").Copy()
$ws.Range("A50:H50").PasteSpecial(-4122)
$ws.Rows("50:50").RowHeight = 15.75

$ws.Range("A50").Value = "Dreadlord Shadow"
$ws.Range("B50").Value = "Shadow"
$ws.Range("C50").Value = "Conquering Dreadlord"
$ws.Range("D50").Value = "70OYS"
$ws.Range("E50").Value = "Conquering Dreadlord"
$ws.Range("F50").Value = "70OYS"
$ws.Range("G50").Value = "None"
$ws.Range("H50").Value = "None"

# ---------------------------------------------------------------------------
# 3) Update the window/selection state to reflect where the user ended up.
# ---------------------------------------------------------------------------
$ws.Range("G52").Select()
